# Readme_IBIS_IP_2.2.pptx — "Add files via upload"
#
# 1) Drop the extra trailing slide (old slide 6 / sldId 282).
# 2) Retitle the cover slide from "Version 2.3" to "Version 2.2".
# 3) On the four "Structure of include-relations" slides, collapse the
#    "...IBIS-IP " / "2.3 " / "<en-dash> V x.y Services " run triplet
#    into a single run reading "...IBIS-IP 2.2 <en-dash> V x.y Services ".

$p = $ppt.ActivePresentation

# --- 1) remove the obsolete last slide -------------------------------------
$p.Slides.Item($p.Slides.Count).Delete()

# --- 2) cover slide version string -----------------------------------------
$cover = $p.Slides.Item(1)
$title = $cover.Shapes.Title
$tr = $title.TextFrame.TextRange
# "Version " (8 chars) stays as-is; "2.3" (the 2nd run) is replaced/merged.
$tail = $tr.Characters(9, $tr.Length - 8)
$tail.Text = ""
$head = $tr.Characters(1, $tr.Length)
$head.Text = "Version 2.2"

# --- 3) the four "...IBIS-IP 2.3 - V x.y Services" title slides ------------
$dash = [char]8211
$versions = @{ 2 = "V 1.0"; 3 = "V 2.0"; 4 = "V 2.1"; 5 = "V 2.2" }

foreach ($idx in $versions.Keys) {
    $slide = $p.Slides.Item($idx)
    $ver = $versions[$idx]

    $t = $slide.Shapes.Title
    $tr2 = $t.TextFrame.TextRange

    # Fixed prefix "Structure of include" (20 chars) + "-relations in IBIS-IP "
    # (22 chars) = run ends at char 42; "2.3 " occupies 43-46; the final run
    # ("<dash> V x.y Services ") starts at 47 and runs to the end.
    $tail2 = $tr2.Characters(43, $tr2.Length - 42)
    $tail2.Text = ""

    $run1 = $tr2.Characters(21, $tr2.Length - 20)
    $run1.Text = "-relations in IBIS-IP 2.2 " + $dash + " " + $ver + " Services "
}
